$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.760.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.784.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '662.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.86%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.37%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.782.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.25%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.31%  '

# Row 10
$ws.Range("E10").Value = '  -0.71%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.458'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.08%  '

# Row 12
$ws.Range("E12").Value = '  +4.27%  '

# Row 13
$ws.Range("E13").Value = '  -3.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.32%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.425.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.786.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.693.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.39%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.114'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '470.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.57%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("E23").Value = '  +1.01%  '

# Row 24
$ws.Range("E24").Value = '  -3.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.80%  '

# Row 26
$ws.Range("E26").Value = '  +1.76%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.65%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.64%  '

# Row 29
$ws.Range("E29").Value = '  +0.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.933.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.83%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.77%  '

# Row 32
$ws.Range("E32").Value = '  +2.98%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.66%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.71%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +15.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.740.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
$ws.Range("E39").Value = '  -1.43%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$ws.Range("E41").Value = '  -0.06%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.961'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.49%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.22%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.32%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.43%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.301'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.71%  '

# Row 50
$ws.Range("E50").Value = '  +1.16%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.57%  '
